# Updates the "cryptos" price/volume snapshot table (rows 2-51) to match the
# latest scrape. Price values are stored as plain text (not numbers) in the
# source data, so for any replacement price that parses as a clean decimal
# number we force the cell's number format to Text ("@") first -- otherwise
# the COM layer would silently coerce it to a float (e.g. "211.39" ->
# 211.38999999999999, or "0.0847" -> scientific notation) and/or drop the
# leading/trailing zeros the source keeps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.695.08'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.596.35'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.39'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.67'
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '1.823.21'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').Value = '1.596.93'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.92'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '26.669.41'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('E18').Value = '  -2.31%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.01'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '209.09'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.10'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.99'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.64'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.13'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.32'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0507'
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.26'
$ws.Range('E34').Value = '  +17.82%  '
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('D36').Value = '1.271.27'
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.48'
$ws.Range('E37').Value = '  -0.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.598'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.777'
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.16'
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.58'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').Value = '1.734.13'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.32'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('E47').Value = '  -0.60%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.103'
$ws.Range('E48').Value = '  +2.48%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0103'
$ws.Range('E49').Value = '  -2.76%  '
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('E51').Value = '  +2.19%  '
